# Apply updated cryptocurrency price / volume(1h) figures to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.185.96'
$ws.Range('E2').Value = '  -2.04%  '
$ws.Range('D3').Value = '2.509.89'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '571.58'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '166.30'
$ws.Range('E6').Value = '  -1.99%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +1.39%  '
$ws.Range('D9').Value = '2.507.39'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.160'
$ws.Range('E10').Value = '  -2.15%  '
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('E12').Value = '  +4.06%  '
$ws.Range('E13').Value = '  +2.53%  '
$ws.Range('D14').Value = '2.976.95'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').Value = '69.143.98'
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('E16').Value = '  -2.84%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '24.84'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = '2.515.25'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('E19').Value = '  -1.40%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.68'
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '349.08'
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.93'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.96'
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '69.98'
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('E26').Value = '  -1.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.92'
$ws.Range('E27').Value = '  -3.14%  '
$ws.Range('D28').Value = '2.648.23'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('E30').Value = '  -1.93%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.86'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '462.46'
$ws.Range('E32').Value = '  -3.70%  '
$ws.Range('E33').Value = '  -2.57%  '
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('E36').Value = '  +0.78%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '157.89'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.98'
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.52'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.75'
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.318'
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('E43').Value = '  -2.85%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '38.29'
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.13'
$ws.Range('E45').Value = '  -13.26%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.26'
$ws.Range('E46').Value = '  -6.11%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '141.58'
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.526'
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('E49').Value = '  -1.33%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0728'
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.55'
$ws.Range('E51').Value = '  -3.73%  '
